{"js": "// Replace each old cell/date text with its corresponding new text.\n// The mapping below is derived from the unified OOXML diff: each entry is\n// a unique \"old text\" -> \"new text\" pair (one w:t run each), so a plain\n// case-sensitive exact-text search-and-replace for every pair reproduces\n// the edit deterministically.\nconst replacements = [\n  [\"2025-02-07 Friday\", \"2025-02-08 Saturday\"],\n  [\"283\u00f72=141, 1\", \"404\u00f76=67, 2\"],\n  [\"576\u00f72=288, 0\", \"548\u00f73=182, 2\"],\n  [\"178\u00f72=89, 0\", \"303\u00f76=50, 3\"],\n  [\"447\u00f72=223, 1\", \"921\u00f79=102, 3\"],\n  [\"215\u00f73=71, 2\", \"712\u00f77=101, 5\"],\n  [\"909\u00f79=101, 0\", \"908\u00f78=113, 4\"],\n  [\"191\u00f78=23, 7\", \"777\u00f73=259, 0\"],\n  [\"483\u00f76=80, 3\", \"610\u00f77=87, 1\"],\n  [\"442\u00f74=110, 2\", \"408\u00f72=204, 0\"],\n  [\"423\u00f72=211, 1\", \"532\u00f73=177, 1\"],\n  [\"549\u00f72=274, 1\", \"880\u00f74=220, 0\"],\n  [\"104\u00f76=17, 2\", \"652\u00f76=108, 4\"],\n  [\"194\u00f77=27, 5\", \"695\u00f72=347, 1\"],\n  [\"270\u00f78=33, 6\", \"614\u00f77=87, 5\"],\n  [\"956\u00f72=478, 0\", \"469\u00f72=234, 1\"],\n  [\"593\u00f72=296, 1\", \"837\u00f78=104, 5\"],\n  [\"216\u00f77=30, 6\", \"422\u00f73=140, 2\"],\n  [\"262\u00f74=65, 2\", \"448\u00f72=224, 0\"],\n  [\"935\u00f76=155, 5\", \"805\u00f76=134, 1\"],\n  [\"539\u00f78=67, 3\", \"788\u00f79=87, 5\"],\n  [\"374\u00f73=124, 2\", \"971\u00f77=138, 5\"],\n  [\"293\u00f78=36, 5\", \"672\u00f77=96, 0\"],\n  [\"811\u00f73=270, 1\", \"889\u00f75=177, 4\"],\n  [\"183\u00f72=91, 1\", \"395\u00f78=49, 3\"],\n  [\"753\u00f75=150, 3\", \"649\u00f78=81, 1\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  // Replace every occurrence found (expected to be exactly one, since the\n  // old values are unique within the document).\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each date/exercise text with its corresponding new text.\n# The mapping below is derived from the unified OOXML diff: every entry is a\n# unique \"old text\" -> \"new text\" pair (one w:t run each in the source doc),\n# so a plain Find/Replace (wdReplaceAll restricted to an exact match) for each\n# pair reproduces the edit deterministically.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-02-07 Friday\", \"2025-02-08 Saturday\"),\n    @(\"283\u00f72=141, 1\", \"404\u00f76=67, 2\"),\n    @(\"576\u00f72=288, 0\", \"548\u00f73=182, 2\"),\n    @(\"178\u00f72=89, 0\", \"303\u00f76=50, 3\"),\n    @(\"447\u00f72=223, 1\", \"921\u00f79=102, 3\"),\n    @(\"215\u00f73=71, 2\", \"712\u00f77=101, 5\"),\n    @(\"909\u00f79=101, 0\", \"908\u00f78=113, 4\"),\n    @(\"191\u00f78=23, 7\", \"777\u00f73=259, 0\"),\n    @(\"483\u00f76=80, 3\", \"610\u00f77=87, 1\"),\n    @(\"442\u00f74=110, 2\", \"408\u00f72=204, 0\"),\n    @(\"423\u00f72=211, 1\", \"532\u00f73=177, 1\"),\n    @(\"549\u00f72=274, 1\", \"880\u00f74=220, 0\"),\n    @(\"104\u00f76=17, 2\", \"652\u00f76=108, 4\"),\n    @(\"194\u00f77=27, 5\", \"695\u00f72=347, 1\"),\n    @(\"270\u00f78=33, 6\", \"614\u00f77=87, 5\"),\n    @(\"956\u00f72=478, 0\", \"469\u00f72=234, 1\"),\n    @(\"593\u00f72=296, 1\", \"837\u00f78=104, 5\"),\n    @(\"216\u00f77=30, 6\", \"422\u00f73=140, 2\"),\n    @(\"262\u00f74=65, 2\", \"448\u00f72=224, 0\"),\n    @(\"935\u00f76=155, 5\", \"805\u00f76=134, 1\"),\n    @(\"539\u00f78=67, 3\", \"788\u00f79=87, 5\"),\n    @(\"374\u00f73=124, 2\", \"971\u00f77=138, 5\"),\n    @(\"293\u00f78=36, 5\", \"672\u00f77=96, 0\"),\n    @(\"811\u00f73=270, 1\", \"889\u00f75=177, 4\"),\n    @(\"183\u00f72=91, 1\", \"395\u00f78=49, 3\"),\n    @(\"753\u00f75=150, 3\", \"649\u00f78=81, 1\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    # 2 = wdReplaceAll\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\nWrite-Output \"done\"\n"}
